$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B6").Value = "draft"
$ws.Range("B8").Value = "2023-08-01T16:12:28+00:00"
